# Week 15 logged + Week 16 simulated: append new per-game samples to the
# running space-separated series on YDS/ST, and update the season-to-date
# totals on OFF/2minOFF/DEF/2minDEF/ST/TURNS/PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this week's rush/pass yardage-per-play samples.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 1 9 7 13 -1 1 9 3 4 3 3 -2 16 -1 8 0 4 4 11 -1 2 0 1 -7 5 -2 23 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 12 27 6 8 9 0 11 6 2 9 11 8 12 11 23 11 6 13 3 16 10 10 17 19"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 7 1 4 7 10 2 3 10 9 2 0 -1 8 0 5 19 11 -2 6 2 6 0 3 1 3 3 5 3 2 2 4 2 1"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 5 3 0 2 18 12 4 9 7 7 7 13"

# ---------------------------------------------------------------------
# OFF sheet: season totals.
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 182
$offWs.Range("E2").Value = 14
$offWs.Range("F2").Value = 68
$offWs.Range("G2").Value = 36
$offWs.Range("J2").Value = 26
$offWs.Range("L2").Value = 209
$offWs.Range("M2").Value = 129
$offWs.Range("O2").Value = 24
$offWs.Range("P2").Value = 11
$offWs.Range("Q2").Value = 444

$offWs.Range("B3").Value = 9
$offWs.Range("C3").Value = 136
$offWs.Range("E3").Value = 34
$offWs.Range("F3").Value = 69
$offWs.Range("G3").Value = 33
$offWs.Range("H3").Value = 31
$offWs.Range("I3").Value = 47
$offWs.Range("J3").Value = 39
$offWs.Range("N3").Value = 21

# ---------------------------------------------------------------------
# DEF sheet: season totals.
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("B2").Value = 4
$defWs.Range("C2").Value = 173
$defWs.Range("E2").Value = 15
$defWs.Range("F2").Value = 53
$defWs.Range("G2").Value = 52
$defWs.Range("H2").Value = 8
$defWs.Range("J2").Value = 29
$defWs.Range("L2").Value = 192
$defWs.Range("M2").Value = 126
$defWs.Range("O2").Value = 22
$defWs.Range("Q2").Value = 419

$defWs.Range("C3").Value = 134
$defWs.Range("E3").Value = 23
$defWs.Range("F3").Value = 77
$defWs.Range("G3").Value = 31
$defWs.Range("I3").Value = 38
$defWs.Range("J3").Value = 48
$defWs.Range("N3").Value = 22

# ---------------------------------------------------------------------
# ST sheet: season totals + appended per-game samples.
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 53
$stWs.Range("D2").Value = 51
$stWs.Range("J2").Value = 112
$stWs.Range("K2").Value = 105
$stWs.Range("L2").Value = 56

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 61 66"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 18 33"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 34 27"
$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 44 45 72"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 11 3 2"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 6 -1 0 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: season totals.
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("C2").Value = 5
$turnsWs.Range("D2").Value = 9
$turnsWs.Range("E2").Value = 8

$turnsWs.Range("D3").Value = 6
$turnsWs.Range("E3").Value = 7

# ---------------------------------------------------------------------
# PEN sheet: season totals.
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("B2").Value = 18
$penWs.Range("D2").Value = 7
$penWs.Range("B3").Value = 12
